$d = $word.ActiveDocument

# 1) Merge the split "Reviewed "/"Github"/" basics" runs (with spell-check
#    proofErr markers) back into a single run with no xml:space artifacts.
$d.Content.Find.Execute("Reviewed Github basics", $false, $false, $false, `
    $false, $false, $true, 1, $false, "Reviewed Github basics", 2) | Out-Null

# 2) Same for "Showed "/"Github"/" example".
$d.Content.Find.Execute("Showed Github example", $false, $false, $false, `
    $false, $false, $true, 1, $false, "Showed Github example", 2) | Out-Null

# 3) Append one new paragraph at the very end of the document, right after
#    the existing trailing empty paragraph, holding
#    "Testing GitHub - Chris Nevares".
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$runProps = '<w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>'
$enDash = [char]0x2013

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter() | Out-Null

$textPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$textXml = "<w:p $wNs><w:pPr>$runProps</w:pPr><w:r>$runProps<w:t>Testing GitHub $enDash Chris Nevares</w:t></w:r></w:p>"
$textPara.Range.InsertXML($textXml) | Out-Null

Write-Host "Final paragraph count:" $d.Paragraphs.Count
